$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(7, 8).Value = 1500
$ws.Cells.Item(7, 9).Value = 1500
$ws.Cells.Item(7, 10).Value = 0
$ws.Cells.Item(7, 11).Value = 1500
$ws.Cells.Item(7, 12).Value = 0
$ws.Cells.Item(7, 13).Value = -1388
$ws.Cells.Item(7, 14).ClearContents()

$ws.Cells.Item(10, 8).Value = 17417.334
$ws.Cells.Item(10, 9).Value = 12501.333
$ws.Cells.Item(10, 11).Value = 12501.333
$ws.Cells.Item(10, 13).Value = -12208.333

$ws.Cells.Item(14, 8).Value = 1500
$ws.Cells.Item(14, 9).Value = 1500
$ws.Cells.Item(14, 10).Value = 0
$ws.Cells.Item(14, 11).Value = 1500
$ws.Cells.Item(14, 12).Value = 0
$ws.Cells.Item(14, 13).Value = -1309
$ws.Cells.Item(14, 14).ClearContents()

$ws.Cells.Item(16, 8).Value = 1500
$ws.Cells.Item(16, 9).Value = 1500
$ws.Cells.Item(16, 11).Value = 1500
$ws.Cells.Item(16, 13).Value = -1270

$ws.Cells.Item(38, 8).Value = 967.5
$ws.Cells.Item(38, 9).Value = 140.72728
$ws.Cells.Item(38, 11).Value = 422.18184
$ws.Cells.Item(38, 13).Value = -50.18184000000002

$ws.Cells.Item(58, 8).Value = 1721.5
$ws.Cells.Item(58, 10).Value = 3000
$ws.Cells.Item(58, 12).Value = 9000
$ws.Cells.Item(58, 14).Value = -9300

$ws.Cells.Item(98, 8).Value = 943.28125
$ws.Cells.Item(98, 9).Value = 943.28125
$ws.Cells.Item(98, 11).Value = 943.28125
$ws.Cells.Item(98, 13).Value = 554.71875

$ws.Cells.Item(103, 8).Value = 575.0612
$ws.Cells.Item(103, 9).Value = 495.2
$ws.Cells.Item(103, 10).Value = 930
$ws.Cells.Item(103, 11).Value = 1485.6
$ws.Cells.Item(103, 12).Value = 2790
$ws.Cells.Item(103, 13).Value = -899.5999999999999
$ws.Cells.Item(103, 14).Value = -3962

$ws.Cells.Item(110, 8).Value = 44333
$ws.Cells.Item(110, 10).Value = 44333
$ws.Cells.Item(110, 12).Value = 44333
$ws.Cells.Item(110, 14).Value = -52513

$ws.Cells.Item(122, 8).Value = 943.28125
$ws.Cells.Item(122, 9).Value = 943.28125
$ws.Cells.Item(122, 11).Value = 2829.84375
$ws.Cells.Item(122, 13).Value = -379.84375

$ws.Cells.Item(129, 8).Value = 26424.125
$ws.Cells.Item(129, 9).Value = 20278.8
$ws.Cells.Item(129, 11).Value = 60836.39999999999
$ws.Cells.Item(129, 13).Value = -55836.39999999999

$ws.Cells.Item(132, 8).Value = 4852.953
$ws.Cells.Item(132, 9).Value = 5040.983
$ws.Cells.Item(132, 11).Value = 15122.949
$ws.Cells.Item(132, 13).Value = -12592.949

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(88, 8).Value = 1354.762
$ws.Cells.Item(88, 10).Value = 1557.0667
$ws.Cells.Item(88, 12).Value = 1557.0667
$ws.Cells.Item(88, 14).Value = -2369.0667

$ws.Cells.Item(91, 8).Value = 1354.762
$ws.Cells.Item(91, 10).Value = 1557.0667
$ws.Cells.Item(91, 12).Value = 1557.0667
$ws.Cells.Item(91, 14).Value = -4365.0667

$ws.Cells.Item(132, 8).Value = 596831.3
$ws.Cells.Item(132, 9).Value = 695847.6
$ws.Cells.Item(132, 10).Value = 2733.3333
$ws.Cells.Item(132, 11).Value = 2087542.8
$ws.Cells.Item(132, 12).Value = 8199.999899999999
$ws.Cells.Item(132, 13).Value = -2085012.8
$ws.Cells.Item(132, 14).Value = -13259.9999

$ws.Cells.Item(135, 8).Value = 0
$ws.Cells.Item(135, 10).Value = 0
$ws.Cells.Item(135, 12).Value = 0
$ws.Cells.Item(135, 14).ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80, 8).Value = 261.625
$ws.Cells.Item(80, 10).Value = 168.25
$ws.Cells.Item(80, 12).Value = 168.25
$ws.Cells.Item(80, 14).Value = -2164.25

$ws.Cells.Item(83, 8).Value = 261.625
$ws.Cells.Item(83, 10).Value = 168.25
$ws.Cells.Item(83, 12).Value = 841.25
$ws.Cells.Item(83, 14).Value = -10825.25

$ws.Cells.Item(86, 8).Value = 2694.353
$ws.Cells.Item(86, 9).Value = 2836.818
$ws.Cells.Item(86, 10).Value = 2433.1667
$ws.Cells.Item(86, 11).Value = 2836.818
$ws.Cells.Item(86, 12).Value = 2433.1667
$ws.Cells.Item(86, 13).Value = -1713.818
$ws.Cells.Item(86, 14).Value = -4679.1667

$ws.Cells.Item(89, 8).Value = 2694.353
$ws.Cells.Item(89, 9).Value = 2836.818
$ws.Cells.Item(89, 10).Value = 2433.1667
$ws.Cells.Item(89, 11).Value = 14184.09
$ws.Cells.Item(89, 12).Value = 12165.8335
$ws.Cells.Item(89, 13).Value = -8568.09
$ws.Cells.Item(89, 14).Value = -23397.8335

$ws.Cells.Item(99, 8).Value = 5806.375
$ws.Cells.Item(99, 9).Value = 6006.696
$ws.Cells.Item(99, 11).Value = 6006.696
$ws.Cells.Item(99, 13).Value = -4508.696

$ws.Cells.Item(107, 8).Value = 11906273
$ws.Cells.Item(107, 9).Value = 13159361
$ws.Cells.Item(107, 11).Value = 13159361
$ws.Cells.Item(107, 13).Value = -13157441

$ws.Cells.Item(134, 8).Value = 3972732.2
$ws.Cells.Item(134, 9).Value = 2054.8684
$ws.Cells.Item(134, 11).Value = 6164.6052
$ws.Cells.Item(134, 13).Value = -3629.6052

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 2198011.2
$ws.Cells.Item(58, 9).Value = 2321.0454
$ws.Cells.Item(58, 11).Value = 2321.0454
$ws.Cells.Item(58, 13).Value = -2118.0454

$ws.Cells.Item(59, 8).Value = 31000
$ws.Cells.Item(59, 9).Value = 0
$ws.Cells.Item(59, 11).Value = 0
$ws.Cells.Item(59, 13).ClearContents()

$ws.Cells.Item(86, 8).Value = 28999.777
$ws.Cells.Item(86, 9).Value = 31874.75
$ws.Cells.Item(86, 10).Value = 6000
$ws.Cells.Item(86, 11).Value = 31874.75
$ws.Cells.Item(86, 12).Value = 6000
$ws.Cells.Item(86, 13).Value = -30751.75
$ws.Cells.Item(86, 14).Value = -8246

$ws.Cells.Item(89, 8).Value = 28999.777
$ws.Cells.Item(89, 9).Value = 31874.75
$ws.Cells.Item(89, 10).Value = 6000
$ws.Cells.Item(89, 11).Value = 159373.75
$ws.Cells.Item(89, 12).Value = 30000
$ws.Cells.Item(89, 13).Value = -153757.75
$ws.Cells.Item(89, 14).Value = -41232

$ws.Cells.Item(134, 8).Value = 1931.8889
$ws.Cells.Item(134, 9).Value = 1735.875
$ws.Cells.Item(134, 10).Value = 3500
$ws.Cells.Item(134, 11).Value = 5207.625
$ws.Cells.Item(134, 12).Value = 10500
$ws.Cells.Item(134, 13).Value = -2672.625
$ws.Cells.Item(134, 14).Value = -15570

$ws.Cells.Item(136, 8).Value = 2198011.2
$ws.Cells.Item(136, 9).Value = 2321.0454
$ws.Cells.Item(136, 11).Value = 6963.1362
$ws.Cells.Item(136, 13).Value = -4413.1362

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 3593518.2
$ws.Cells.Item(4, 9).Value = 5257055
$ws.Cells.Item(4, 10).Value = 127816.586
$ws.Cells.Item(4, 11).Value = 15771165
$ws.Cells.Item(4, 12).Value = 383449.758
$ws.Cells.Item(4, 13).Value = -15771053
$ws.Cells.Item(4, 14).Value = -383673.758

$ws.Cells.Item(26, 8).Value = 598.6875
$ws.Cells.Item(26, 9).Value = 250
$ws.Cells.Item(26, 10).Value = 947.375
$ws.Cells.Item(26, 11).Value = 750
$ws.Cells.Item(26, 12).Value = 2842.125
$ws.Cells.Item(26, 13).Value = -462
$ws.Cells.Item(26, 14).Value = -3418.125

$ws.Cells.Item(32, 8).Value = 3948.8
$ws.Cells.Item(32, 9).Value = 3463.3333
$ws.Cells.Item(32, 10).Value = 4156.857
$ws.Cells.Item(32, 11).Value = 10389.9999
$ws.Cells.Item(32, 12).Value = 12470.571
$ws.Cells.Item(32, 13).Value = -10106.9999
$ws.Cells.Item(32, 14).Value = -13036.571

$ws.Cells.Item(33, 8).Value = 2126.6
$ws.Cells.Item(33, 9).Value = 116.5
$ws.Cells.Item(33, 10).Value = 3466.6667
$ws.Cells.Item(33, 11).Value = 699
$ws.Cells.Item(33, 12).Value = 20800.0002
$ws.Cells.Item(33, 13).Value = -416
$ws.Cells.Item(33, 14).Value = -21366.0002

$ws.Cells.Item(41, 8).Value = 0
$ws.Cells.Item(41, 9).Value = 0
$ws.Cells.Item(41, 11).Value = 0
$ws.Cells.Item(41, 13).ClearContents()

$ws.Cells.Item(44, 8).Value = 2859.3
$ws.Cells.Item(44, 9).Value = 465.66666
$ws.Cells.Item(44, 10).Value = 6449.75
$ws.Cells.Item(44, 11).Value = 1396.99998
$ws.Cells.Item(44, 12).Value = 19349.25
$ws.Cells.Item(44, 13).Value = -998.9999800000001
$ws.Cells.Item(44, 14).Value = -20145.25

$ws.Cells.Item(75, 8).Value = 1899.6
$ws.Cells.Item(75, 9).Value = 499.5
$ws.Cells.Item(75, 10).Value = 2833
$ws.Cells.Item(75, 11).Value = 1498.5
$ws.Cells.Item(75, 12).Value = 8499
$ws.Cells.Item(75, 13).Value = -500.5
$ws.Cells.Item(75, 14).Value = -10495

$ws.Cells.Item(78, 8).Value = 1899.6
$ws.Cells.Item(78, 9).Value = 499.5
$ws.Cells.Item(78, 10).Value = 2833
$ws.Cells.Item(78, 11).Value = 4495.5
$ws.Cells.Item(78, 12).Value = 25497
$ws.Cells.Item(78, 13).Value = 496.5
$ws.Cells.Item(78, 14).Value = -35481

$ws.Cells.Item(131, 8).Value = 2711621.2
$ws.Cells.Item(131, 10).Value = 3474025.8
$ws.Cells.Item(131, 12).Value = 10422077.4
$ws.Cells.Item(131, 14).Value = -10432157.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 3692.1936
$ws.Cells.Item(122, 9).Value = 3743.5908
$ws.Cells.Item(122, 10).Value = 3566.5557
$ws.Cells.Item(122, 11).Value = 11230.7724
$ws.Cells.Item(122, 12).Value = 10699.6671
$ws.Cells.Item(122, 13).Value = -8780.7724
$ws.Cells.Item(122, 14).Value = -15599.6671

$ws.Cells.Item(132, 8).Value = 9546.370999999999
$ws.Cells.Item(132, 9).Value = 11294
$ws.Cells.Item(132, 11).Value = 33882
$ws.Cells.Item(132, 13).Value = -31352

$ws.Cells.Item(135, 8).Value = 78777.234
$ws.Cells.Item(135, 10).Value = 78777.234
$ws.Cells.Item(135, 12).Value = 78777.234
$ws.Cells.Item(135, 14).Value = -88917.234

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 4075.2285
$ws.Cells.Item(46, 9).Value = 1544.25
$ws.Cells.Item(46, 10).Value = 4401.8066
$ws.Cells.Item(46, 11).Value = 1544.25
$ws.Cells.Item(46, 12).Value = 4401.8066
$ws.Cells.Item(46, 13).Value = -1356.25
$ws.Cells.Item(46, 14).Value = -4777.8066

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 8038.077
$ws.Cells.Item(62, 10).Value = 8600.6
$ws.Cells.Item(62, 12).Value = 8600.6
$ws.Cells.Item(62, 14).Value = -9848.6

$ws.Cells.Item(65, 8).Value = 8038.077
$ws.Cells.Item(65, 10).Value = 8600.6
$ws.Cells.Item(65, 12).Value = 43003
$ws.Cells.Item(65, 14).Value = -49243

$ws.Cells.Item(132, 8).Value = 15154320
$ws.Cells.Item(132, 9).Value = 16669312
$ws.Cells.Item(132, 10).Value = 4400
$ws.Cells.Item(132, 11).Value = 50007936
$ws.Cells.Item(132, 12).Value = 13200
$ws.Cells.Item(132, 13).Value = -50005406
$ws.Cells.Item(132, 14).Value = -18260
